$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'94.401.80"
$ws.Range("E2").Value = "  +3.26%  "
$ws.Range("D3").Value = "'3.135.59"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'241.55"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "'617.40"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +5.66%  "
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'0.832"
$ws.Range("E10").Value = "  +14.47%  "
$ws.Range("D11").Value = "'3.129.48"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "'93.973.88"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'35.02"
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("D16").Value = "'5.42"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "'3.717.37"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "'3.160.47"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("D19").Value = "'3.71"
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").Value = "'15.03"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("D22").Value = "'452.94"
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").Value = "'0.0000203"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").Value = "'9.09"
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").Value = "'8.30"
$ws.Range("E25").Value = "  +5.53%  "
$ws.Range("D26").Value = "'5.71"
$ws.Range("E26").Value = "  +2.13%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'86.80"
$ws.Range("E27").Value = "  +7.71%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'12.24"
$ws.Range("E28").Value = "  +5.80%  "
$ws.Range("D29").Value = "'3.298.73"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.266"
$ws.Range("E30").Value = "  +17.59%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +9.24%  "
$ws.Range("D33").Value = "'0.128"
$ws.Range("E33").Value = "  -5.83%  "
$ws.Range("D34").Value = "'9.39"
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "'8.03"
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").Value = "'0.164"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "'26.33"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").Value = "'1.92"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "'0.458"
$ws.Range("E40").Value = "  +6.73%  "
$ws.Range("D41").Value = "'484.28"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").Value = "'1.29"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'23.94"
$ws.Range("E43").Value = "  +7.94%  "
$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D44").Value = "'3.77"
$ws.Range("E44").Value = "  -8.34%  "
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("D47").Value = "'159.80"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "'0.699"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "'4.49"
$ws.Range("E50").Value = "  +2.94%  "
$ws.Range("D51").Value = "'1.33"
$ws.Range("E51").Value = "  -1.38%  "
